$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) so numeric-looking text (e.g. "212.34", "1.00")
# is written back as text instead of being auto-converted to a number by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.200.19'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.604.59'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '212.34'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.484'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").Value = '0.0612'
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = '18.13'
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.828.97'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.589.25'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("D16").Value = '26.190.76'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '61.70'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = ('0.0{0}0727' -f [char]0x2083)
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '200.48'
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("D22").Value = '9.26'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '6.00'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("D25").Value = '144.17'
$ws.Range("E25").Value = '  +2.05%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("D28").Value = '15.17'
$ws.Range("E29").Value = '  +1.48%  '
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +3.94%  '
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +2.37%  '
$ws.Range("D33").Value = '2.92'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +0.70%  '
$ws.Range("D36").Value = '1.165.55'
$ws.Range("E36").Value = '  +5.06%  '
$ws.Range("E37").Value = '  +4.91%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = '0.785'
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("D41").Value = '0.495'
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("D43").Value = '5.28'
$ws.Range("E43").Value = '  +3.79%  '
$ws.Range("D44").Value = '1.741.02'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").Value = '91.55'
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("D47").Value = '54.01'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = ('0.0{0}0971' -f [char]0x2087)
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.408'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.05%  '

# Restore the original (default) cell style now that the text values are set.
$priceRange.Style = "Normal"
